$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing weekly price records (rows 6-17) to reflect the
# newly inserted observation and the resulting downward shift of the
# remaining weekly records.
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11500
$ws.Range("S6").Value = 767
$ws.Range("D7").Value = 45033
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("Q7").Value = '$/caja 18 kilos empedrada'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 722
$ws.Range("T7").Value = 18
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 667
$ws.Range("D9").Value = 45062
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("S9").Value = 722
$ws.Range("D10").Value = 45062
$ws.Range("Q10").Value = '$/caja 18 kilos empedrada'
$ws.Range("D11").Value = 45021
$ws.Range("L11").Value = 'Primera'
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 667
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("S12").Value = 722
$ws.Range("D13").Value = 45050
$ws.Range("L13").Value = 'Primera'
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 667
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 40
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("S14").Value = 722
$ws.Range("D15").Value = 45043
$ws.Range("M15").Value = 50
$ws.Range("Q15").Value = '$/caja 18 kilos empedrada'
$ws.Range("D16").Value = 45020
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 12000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("S16").Value = 667
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 13000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 13000
$ws.Range("S17").Value = 722
# New row 18 (appended at the end of the weekly series).
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C18").Value = 'Ñuble'
$ws.Range("D18").Value = 45040
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 'Fruta'
$ws.Range("G18").Value = 100104
$ws.Range("H18").Value = 'Frutos de pepita'
$ws.Range("I18").Value = 100104003
$ws.Range("J18").Value = 'Membrillo'
$ws.Range("K18").Value = 'Champion'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 12000
$ws.Range("P18").Value = 12000
$ws.Range("Q18").Value = '$/caja 18 kilos empedrada'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 667
$ws.Range("T18").Value = 18

# Match the date formatting used by the other rows in column D.
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
